# Updates the cryptos price/volume snapshot values (columns D and E) and
# fixes the HuobiToken/ImmutableX row ordering (rows 34-35), matching the
# latest scraped data commit.
#
# Note: several "Price" values (column D) are plain decimal numbers
# (e.g. "261.58"). Assigning those as a plain string would make Excel
# auto-convert the cell to a numeric type, which would not match the
# original inline/shared-string text cells. To keep them as text we set
# them with a leading apostrophe (forces text entry) and then reset the
# cell Style back to "Normal" so no extra number-format/quote-prefix
# style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.461.96"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.837.72"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'261.58"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5386"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3016"
$ws.Range("E8").Value = "  -7.28%  "
$ws.Range("D9").Value = "'0.06881"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "'17.68"
$ws.Range("E10").Value = "  -6.92%  "
$ws.Range("D11").Value = "'0.7380"
$ws.Range("E11").Value = "  -5.90%  "
$ws.Range("D12").Value = "1.837.59"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'0.07159"
$ws.Range("E13").Value = "  -8.16%  "
$ws.Range("D14").Value = "'89.20"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'4.993"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'0.000007880"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").Value = "26.490.66"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "2.081.26"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'5.980"
$ws.Range("D24").Value = "'9.207"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "'142.92"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'1.695"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'110.67"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'4.237"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'0.08823"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'4.033"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "'0.04819"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7299"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.919"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "'3.093"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "'2.262"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "'0.01714"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("D40").Value = "'0.4720"
$ws.Range("E40").Value = "  -3.28%  "
$ws.Range("D41").Value = "'0.9053"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'107.94"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("D43").Value = "'5.891"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").Value = "'0.1244"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "'9.005"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'0.4073"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").Value = "'34.82"
$ws.Range("D50").Value = "'0.05774"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "'0.8930"
$ws.Range("E51").Value = "  +0.34%  "

# Reset style for cells that required a quote-prefix to stay text,
# so no extra style index is attached to them (matches original formatting).
$numericTextCells = @("D5","D7","D8","D9","D10","D11","D13","D14","D15","D19","D23","D24","D25","D27","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D46","D47","D48","D49","D50","D51")
foreach ($cell in $numericTextCells) {
    $ws.Range($cell).Style = "Normal"
}
